# Minor update on which columns to show in table 1 and supptable1
# - Add PMID (31142855) into column B, row 4 on both sheets
# - Make "table1_metadata" (sheet1) the active/selected sheet,
#   with B4 selected on both sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("table1_metadata")
$ws2 = $wb.Worksheets.Item("protocal")

# Insert the new PMID value in column B row 4 on both sheets
$ws1.Range("B4").Value = 31142855
$ws2.Range("B4").Value = 31142855

# Update selections on each sheet to B4
$ws1.Range("B4").Select()
$ws2.Range("B4").Select()

# Make table1_metadata the active (displayed/tab-selected) sheet
$ws1.Activate()
$ws1.Range("B4").Select()
